$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 628 entirely ("「きみは何をえらぶ？」" entry), shifting all
# subsequent rows up by one.
$ws.Rows(628).Delete()
